$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1741
$ws.Cells.Item(3, 6).Value = 10173
$ws.Cells.Item(6, 6).Value = 599
$ws.Cells.Item(7, 6).Value = 75
$ws.Cells.Item(8, 6).Value = 1626
$ws.Cells.Item(9, 6).Value = 180
$ws.Cells.Item(10, 6).Value = 384
$ws.Cells.Item(12, 6).Value = 210
$ws.Cells.Item(14, 6).Value = 481
$ws.Cells.Item(15, 6).Value = 1185
$ws.Cells.Item(19, 6).Value = 94
$ws.Cells.Item(20, 6).Value = 357
$ws.Cells.Item(21, 6).Value = 15
$ws.Cells.Item(22, 6).Value = 329
$ws.Cells.Item(23, 6).Value = 103
$ws.Cells.Item(24, 6).Value = 1159
$ws.Cells.Item(25, 6).Value = 701
$ws.Cells.Item(26, 6).Value = 22
$ws.Cells.Item(29, 6).Value = 231
$ws.Cells.Item(31, 6).Value = 455
$ws.Cells.Item(33, 6).Value = 373
$ws.Cells.Item(35, 6).Value = 638
$ws.Cells.Item(36, 6).Value = 747
$ws.Cells.Item(39, 6).Value = 814
$ws.Cells.Item(40, 6).Value = 383
$ws.Cells.Item(41, 6).Value = 333
$ws.Cells.Item(42, 6).Value = 15
$ws.Cells.Item(43, 6).Value = 356
$ws.Cells.Item(44, 6).Value = 78
$ws.Cells.Item(45, 6).Value = 353
$ws.Cells.Item(46, 6).Value = 78

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 89
$ws.Cells.Item(14, 6).Value = 101
$ws.Cells.Item(18, 6).Value = 1097
$ws.Cells.Item(20, 6).Value = 693
$ws.Cells.Item(21, 6).Value = 1102
$ws.Cells.Item(22, 6).Value = 327
$ws.Cells.Item(23, 6).Value = 686
$ws.Cells.Item(24, 6).Value = 76
$ws.Cells.Item(25, 6).Value = 10
$ws.Cells.Item(31, 6).Value = 208
$ws.Cells.Item(34, 6).Value = 157
$ws.Cells.Item(35, 6).Value = 188
$ws.Cells.Item(41, 6).Value = 5
$ws.Cells.Item(42, 6).Value = 5
$ws.Cells.Item(43, 6).Value = 69
$ws.Cells.Item(44, 6).Value = 41

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 807
$ws.Cells.Item(5, 6).Value = 194
$ws.Cells.Item(6, 6).Value = 2525
$ws.Cells.Item(7, 6).Value = 4086
$ws.Cells.Item(8, 6).Value = 64
$ws.Cells.Item(10, 6).Value = 306
$ws.Cells.Item(11, 6).Value = 193

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1741
$ws.Cells.Item(3, 6).Value = 807
$ws.Cells.Item(4, 6).Value = 10173
$ws.Cells.Item(5, 6).Value = 194
$ws.Cells.Item(7, 6).Value = 4086
$ws.Cells.Item(8, 6).Value = 64
$ws.Cells.Item(9, 6).Value = 306
$ws.Cells.Item(10, 6).Value = 306
$ws.Cells.Item(11, 6).Value = 599
$ws.Cells.Item(12, 6).Value = 1626
$ws.Cells.Item(13, 6).Value = 180
$ws.Cells.Item(14, 6).Value = 384
$ws.Cells.Item(16, 6).Value = 210
$ws.Cells.Item(19, 6).Value = 1185
$ws.Cells.Item(23, 6).Value = 101
$ws.Cells.Item(25, 6).Value = 94
$ws.Cells.Item(26, 6).Value = 1097
$ws.Cells.Item(27, 6).Value = 357
$ws.Cells.Item(28, 6).Value = 329
$ws.Cells.Item(29, 6).Value = 1102
$ws.Cells.Item(30, 6).Value = 1159
$ws.Cells.Item(31, 6).Value = 701
$ws.Cells.Item(32, 6).Value = 76
$ws.Cells.Item(34, 6).Value = 231
$ws.Cells.Item(36, 6).Value = 455
$ws.Cells.Item(38, 6).Value = 373
$ws.Cells.Item(40, 6).Value = 638
$ws.Cells.Item(41, 6).Value = 208
$ws.Cells.Item(42, 6).Value = 747
$ws.Cells.Item(44, 6).Value = 814
$ws.Cells.Item(45, 6).Value = 383
$ws.Cells.Item(46, 6).Value = 333
$ws.Cells.Item(48, 6).Value = 356
$ws.Cells.Item(49, 6).Value = 353
$ws.Cells.Item(50, 6).Value = 69
